$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Bug fix: the "counselor edit appointment" endpoint was missing the
# "teacher/" path segment (it should match the sibling
# "appointment/teacher/cancelAppointment" endpoint just below it).
# The document's auto-managed "_GoBack" bookmark (which marks the
# location of the most recent edit) also needs to move from the old
# edit spot (in the cancelAppointment paragraph) to the new edit spot
# (in the editAppointment paragraph).
# -----------------------------------------------------------------

# 1) Remove the existing "_GoBack" bookmark (it currently sits inside
#    the "appointment/teacher/cancelAppointment" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Find the red "POST appointment/editAppointment" paragraph (the
#    counselor-facing endpoint) and locate "editAppointment" inside it.
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if (($t -eq "POST appointment/editAppointment`r") -and ($p.Range.Font.Color -eq 255)) {
        $target = $p
        break
    }
}

$fr = $target.Range
$f = $fr.Find
$f.Text = "editAppointment"
$f.Forward = $true
$f.Wrap = 0
$f.Execute() | Out-Null
$insertPos = $fr.Start

# 3) Re-create the "_GoBack" bookmark right before "editAppointment",
#    i.e. between the new "teacher/" segment and "editAppointment".
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 4) Insert the missing "teacher/" path segment right before
#    "editAppointment" (and before the bookmark we just placed).
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter("teacher/")
